$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data from BFPIaE")

# Updated source data (new "BFPIaE" figures) for Production/Imports/Exports and
# the computed Domestic Use column. Setting .Value on the previously-formula
# cells in column E replaces the shared formula with the pasted static result,
# matching the authored change (values pasted from the updated BFPIaE sheet).
$ws.Range("B5").Value = 13727856983750000
$ws.Range("C5").Value = 116791410000000
$ws.Range("D5").Value = 1888127795000000
$ws.Range("E5").Value = 11956520598750000
$ws.Range("B6").Value = 35698853000000000
$ws.Range("C6").Value = 2791070000000000
$ws.Range("D6").Value = 4483330999999999.5
$ws.Range("E6").Value = 34006591999999996
$ws.Range("B7").Value = 540000000000000
$ws.Range("E7").Value = 7740000000000000
$ws.Range("B11").Value = 123408653368000
$ws.Range("E11").Value = 40800652766378.289
$ws.Range("B12").Value = 18365288407359000
$ws.Range("C12").Value = 59545742184000
$ws.Range("D12").Value = 1381937220783000
$ws.Range("E12").Value = 17042896928760000
$ws.Range("B13").Value = 10682345175000000
$ws.Range("C13").Value = 320229375000000
$ws.Range("D13").Value = 2936702875000000
$ws.Range("E13").Value = 8065871675000000
$ws.Range("B14").Value = 1515620096655000
$ws.Range("C14").Value = 7285809312000
$ws.Range("D14").Value = 132183114996000
$ws.Range("E14").Value = 1390722790971000
$ws.Range("B15").Value = 203604487000000
$ws.Range("C15").Value = 74398997000000
$ws.Range("D15").Value = 11939852000000
$ws.Range("E15").Value = 266063632000000
$ws.Range("B16").Value = 3537592380000000
$ws.Range("C16").Value = 338510340000000
$ws.Range("D16").Value = 393656760000000
$ws.Range("E16").Value = 3482445960000000
$ws.Range("B19").Value = 906213062527442.13
$ws.Range("E19").Value = 906213062527442.13
$ws.Range("B20").Value = 19460537227008000
$ws.Range("C20").Value = 16583078106860000
$ws.Range("D20").Value = 2408884127644000
$ws.Range("E20").Value = 33634731206224000
$ws.Range("B21").Value = 979835237000000
$ws.Range("C21").Value = 433897305000000
$ws.Range("D21").Value = 705652880000000
$ws.Range("E21").Value = 708079662000000
$ws.Range("B22").Value = 3148621106400000
$ws.Range("C22").Value = 222079680900000
$ws.Range("D22").Value = 1475922862260000
$ws.Range("E22").Value = 1894777925040000
$ws.Range("B23").Value = 3564295858911020.5
$ws.Range("E23").Value = 3564295858911020.5
$ws.Range("E24").Value = 8746500000000000

# Re-select the sheet that was active when the file was saved, and select the
# cell that was selected on "Data from BFPIaE" at save time.
$ws.Range("H28").Select()

$ws3 = $wb.Worksheets.Item("PoFDCtAE")
$ws3.Activate()
